# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E3) and "Correspond Handback
# DateTime" (H3) for the f66ea9a6-... file row on both the zh-cn and
# de-de language sheets with the freshly generated handback timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-31 07:09:52"
$zhcn.Range("H3").Value = "2016-03-31 07:10:42"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-31 07:10:09"
$dede.Range("H3").Value = "2016-03-31 07:11:01"
